# Apply the edits captured by the commit "Add files via upload":
#   - B3 (Inventaires mensuels réalisés for "Mai 2025") changes from 3483 to 4449
#   - The active/selected cell on the sheet moves from D4 to C7

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the realised-inventory figure for May 2025
$ws.Range("B3").Value = 4449

# Move the selection to C7 to match the saved cursor position
$ws.Range("C7").Select()
